$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dataset was recomputed after changing which participants'
# decisions count toward the ranking ("include no rank decision in binary").
# This reshuffled a handful of worker rows within each realeffort group
# (their prolificid/name/race moved to a different row) and recalculated the
# re_rank score for every row.

# --- Swap row pairs that exchanged identity (prolificid/name/race) ---
# Rows 10 <-> 13 and 11 <-> 12 (female group)
$c10 = $ws.Range("C10").Value2
$c11 = $ws.Range("C11").Value2
$c12 = $ws.Range("C12").Value2
$c13 = $ws.Range("C13").Value2
$d10 = $ws.Range("D10").Value2
$d11 = $ws.Range("D11").Value2
$d12 = $ws.Range("D12").Value2
$d13 = $ws.Range("D13").Value2
$e10 = $ws.Range("E10").Value2
$e11 = $ws.Range("E11").Value2
$e12 = $ws.Range("E12").Value2
$e13 = $ws.Range("E13").Value2
$h10 = $ws.Range("H10").Value2
$h11 = $ws.Range("H11").Value2
$h12 = $ws.Range("H12").Value2
$h13 = $ws.Range("H13").Value2

$ws.Range("C10").Value2 = $c13
$ws.Range("D10").Value2 = $d13
$ws.Range("E10").Value2 = $e13
$ws.Range("H10").Value2 = $h13

$ws.Range("C11").Value2 = $c12
$ws.Range("D11").Value2 = $d12
$ws.Range("E11").Value2 = $e12
$ws.Range("H11").Value2 = $h12

$ws.Range("C12").Value2 = $c10
$ws.Range("D12").Value2 = $d10
$ws.Range("E12").Value2 = $e10
$ws.Range("H12").Value2 = $h10

$ws.Range("C13").Value2 = $c11
$ws.Range("D13").Value2 = $d11
$ws.Range("E13").Value2 = $e11
$ws.Range("H13").Value2 = $h11

# Rows 16 <-> 17 (male group)
$c16 = $ws.Range("C16").Value2
$c17 = $ws.Range("C17").Value2
$d16 = $ws.Range("D16").Value2
$d17 = $ws.Range("D17").Value2
$e16 = $ws.Range("E16").Value2
$e17 = $ws.Range("E17").Value2
$h16 = $ws.Range("H16").Value2
$h17 = $ws.Range("H17").Value2

$ws.Range("C16").Value2 = $c17
$ws.Range("D16").Value2 = $d17
$ws.Range("E16").Value2 = $e17
$ws.Range("H16").Value2 = $h17

$ws.Range("C17").Value2 = $c16
$ws.Range("D17").Value2 = $d16
$ws.Range("E17").Value2 = $e16
$ws.Range("H17").Value2 = $h16

# Rows 21 <-> 22 (male group)
$c21 = $ws.Range("C21").Value2
$c22 = $ws.Range("C22").Value2
$d21 = $ws.Range("D21").Value2
$d22 = $ws.Range("D22").Value2
$e21 = $ws.Range("E21").Value2
$e22 = $ws.Range("E22").Value2
$h21 = $ws.Range("H21").Value2
$h22 = $ws.Range("H22").Value2

$ws.Range("C21").Value2 = $c22
$ws.Range("D21").Value2 = $d22
$ws.Range("E21").Value2 = $e22
$ws.Range("H21").Value2 = $h22

$ws.Range("C22").Value2 = $c21
$ws.Range("D22").Value2 = $d21
$ws.Range("E22").Value2 = $e21
$ws.Range("H22").Value2 = $h21

# Rows 24 <-> 25 (male group)
$c24 = $ws.Range("C24").Value2
$c25 = $ws.Range("C25").Value2
$d24 = $ws.Range("D24").Value2
$d25 = $ws.Range("D25").Value2
$e24 = $ws.Range("E24").Value2
$e25 = $ws.Range("E25").Value2
$h24 = $ws.Range("H24").Value2
$h25 = $ws.Range("H25").Value2

$ws.Range("C24").Value2 = $c25
$ws.Range("D24").Value2 = $d25
$ws.Range("E24").Value2 = $e25
$ws.Range("H24").Value2 = $h25

$ws.Range("C25").Value2 = $c24
$ws.Range("D25").Value2 = $d24
$ws.Range("E25").Value2 = $e24
$ws.Range("H25").Value2 = $h24

# --- Recomputed re_rank (column G) scores for every row ---
$ws.Range("G2").Value2 = 7.360079283446961
$ws.Range("G3").Value2 = 6.358750456454161
$ws.Range("G4").Value2 = 6.017579664918089
$ws.Range("G5").Value2 = 5.309718579672998
$ws.Range("G6").Value2 = 5.244195657518464
$ws.Range("G7").Value2 = 4.045026469112039
$ws.Range("G8").Value2 = 1.242073243576292
$ws.Range("G9").Value2 = 1.014010395470444
$ws.Range("G10").Value2 = 0.4834459824271087
$ws.Range("G11").Value2 = 0.1753989618967279
$ws.Range("G12").Value2 = 0.1348575153764161
$ws.Range("G13").Value2 = 0.1140016948445168
$ws.Range("G14").Value2 = 13.17322371252606
$ws.Range("G15").Value2 = 8.25975769155853
$ws.Range("G16").Value2 = 7.198709993617562
$ws.Range("G17").Value2 = 7.013420770724821
$ws.Range("G18").Value2 = 5.30118687809812
$ws.Range("G19").Value2 = 5.217927984380697
$ws.Range("G20").Value2 = 5.022934074744907
$ws.Range("G21").Value2 = 4.346124827430741
$ws.Range("G22").Value2 = 4.218944548501164
$ws.Range("G23").Value2 = 3.102299781206951
$ws.Range("G24").Value2 = 2.368353172506046
$ws.Range("G25").Value2 = 2.301699192143767
